# This script re-orders the player roster rows (rows 2-16) in the active
# worksheet so that each (Player, Position, Team) triple lands on its new
# row, while rows 1 (header) and 17-19 remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired state for rows 2 through 16 (columns A, B, C).
$data = @(
    @("Andrew Nembhard", "PG,SG", "Indiana Pacers"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Isaiah Stewart", "PF,C", "Detroit Pistons"),
    @("Nikola Jovic", "PF,C", "Miami Heat"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Cason Wallace", "PG,SG", "Oklahoma City Thunder"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Jaylin Williams", "PF,C", "Oklahoma City Thunder"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
